# update of distances for R6 M II
#
# The "Canon EOS R6 Mark II" block (rows 37-79) gets updated Pixel Distance
# (col G), Scale (col H) and, where applicable, Unit (col I) values. The
# underlying shared-string table also got de-duplicated, but since that has
# no visible effect on cell text we only need to touch the cells whose
# displayed value actually changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => @{ col letter = new value }
$updates = @{
    37 = @{ H = "1000"; I = "µm" }
    38 = @{ G = "1040"; H = "500" }
    39 = @{ G = "836";  H = "200" }
    41 = @{ G = "838";  H = "50"  }
    42 = @{ G = "830";  H = "20"  }
    43 = @{ G = "1068"; H = "4"   }
    44 = @{ G = "656";  H = "2"   }
    45 = @{ G = "830";  H = "2"   }
    46 = @{ G = "1024"; H = "2"   }
    50 = @{ G = "686";  H = "500"; I = "µm" }
    51 = @{ G = "863";  H = "500"; I = "µm" }
    52 = @{ G = "1076"; H = "500"; I = "µm" }
    53 = @{ G = "860";  H = "5"   }
    54 = @{ G = "753";  H = "3"   }
    55 = @{ G = "1005"; H = "3"   }
    56 = @{ G = "850";  H = "2"   }
    57 = @{ G = "1010"; H = "2"   }
    58 = @{ G = "1168"; H = "2"   }
    62 = @{ G = "830";  H = "10"  }
    63 = @{ G = "792";  H = "4"   }
    64 = @{ G = "1112"; H = "4"   }
    65 = @{ G = "975";  H = "3"   }
    66 = @{ G = "864";  H = "2"   }
    67 = @{ G = "974";  H = "2"   }
    68 = @{ G = "1140"; H = "2"   }
    72 = @{ H = "1000"; I = "µm" }
    73 = @{ H = "1000"; I = "µm" }
    74 = @{ H = "1000"; I = "µm" }
    75 = @{ H = "1000"; I = "µm" }
    76 = @{ G = "717";  H = "500"; I = "µm" }
    77 = @{ G = "750";  H = "500"; I = "µm" }
    78 = @{ G = "782";  H = "500"; I = "µm" }
    79 = @{ G = "833";  H = "500"; I = "µm" }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}

# Reflect the author's final cursor position / scroll from the diff.
$ws.Range("G80").Select()
